$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "AI"
$ws.Range("B2").Value = "web development"
$ws.Range("B3").Value = "research"

$ws.Range("B3").Select()
